$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new "2021年" row (row 11) under the existing "2020年" row (row 10).
# Copy row 10's column-A cell first so the new label inherits the same
# bold/centered/bordered style (s="1") used by every other year label.
$ws.Range("A10").Copy($ws.Range("A11"))
$ws.Range("A11").Value = "2021年"

$ws.Range("B11").Value = 743.4400000000001
$ws.Range("C11").Value = 209.19
$ws.Range("D11").Value = 36.27
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = 635.86
$ws.Range("G11").Value = 794.6799999999999
$ws.Range("H11").Value = 38.48
$ws.Range("I11").Value = 1096.76
$ws.Range("J11").Value = 129.14
$ws.Range("K11").Value = 219.64
$ws.Range("L11").Value = 54.85
$ws.Range("M11").Value = 1.83
$ws.Range("N11").Value = 237.56
$ws.Range("O11").Value = 209.87
$ws.Range("P11").Value = 19.69
$ws.Range("Q11").Value = 164.28
$ws.Range("R11").Value = 453.29
$ws.Range("S11").Value = 7.52
$ws.Range("T11").Value = 389.45
$ws.Range("U11").Value = 2.02
$ws.Range("V11").Value = 231.5
$ws.Range("W11").Value = 41.5
$ws.Range("X11").Value = 19.2
$ws.Range("Y11").Value = 1294.95
$ws.Range("Z11").Value = 161.43
$ws.Range("AA11").Value = 134.93
$ws.Range("AB11").Value = 0.06
$ws.Range("AC11").Value = 13194.92
$ws.Range("AD11").Value = 319.33
$ws.Range("AE11").Value = 316.36
$ws.Range("AF11").Value = 815.03
$ws.Range("AG11").Value = 743.88
$ws.Range("AH11").Value = 185.36
$ws.Range("AI11").Value = 293.37
$ws.Range("AJ11").Value = 8.75
$ws.Range("AK11").Value = 616.45
$ws.Range("AL11").Value = 118.19
$ws.Range("AM11").Value = 1521.24
$ws.Range("AN11").Value = 138.63
$ws.Range("AO11").Value = 468.61
$ws.Range("AP11").Value = 283.66
$ws.Range("AQ11").Value = 38.64
